$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 22.05824823270341
$ws.Range("C2").Value = 3.908675865430624
$ws.Range("D2").Value = 4.121503047402403
$ws.Range("F2").Value = 56.17646139172216
$ws.Range("G2").Value = 3.80367133135899
$ws.Range("J2").Value = 10.72874743689601
$ws.Range("K2").Value = 18.05322289120867
$ws.Range("L2").Value = 11.43017458836336
$ws.Range("M2").Value = 19.39557772768245
$ws.Range("N2").Value = 27.24314928207849
$ws.Range("B3").Value = 21.97222296887922
$ws.Range("C3").Value = 3.69404948403072
$ws.Range("D3").Value = 4.128140318954047
$ws.Range("F3").Value = 56.15065750823043
$ws.Range("G3").Value = 3.80676274366199
$ws.Range("J3").Value = 10.74298259935834
$ws.Range("K3").Value = 17.99591887213701
$ws.Range("L3").Value = 11.44843577453275
$ws.Range("M3").Value = 19.40663935041466
$ws.Range("N3").Value = 27.27766621820211
$ws.Range("B4").Value = 21.92440233397227
$ws.Range("C4").Value = 3.555005187740258
$ws.Range("D4").Value = 4.132638628594753
$ws.Range("F4").Value = 56.14365848178091
$ws.Range("G4").Value = 3.808761097760898
$ws.Range("J4").Value = 10.75225109950155
$ws.Range("K4").Value = 17.96475574452864
$ws.Range("L4").Value = 11.4609957131474
$ws.Range("M4").Value = 19.41659691712251
$ws.Range("N4").Value = 27.30070333921554
$ws.Range("B5").Value = 21.9061870528579
$ws.Range("C5").Value = 3.496537203790947
$ws.Range("D5").Value = 4.134578389374596
$ws.Range("F5").Value = 56.14303311443463
$ws.Range("G5").Value = 3.809600730413812
$ws.Range("J5").Value = 10.75616124683982
$ws.Range("K5").Value = 17.95307711021086
$ws.Range("L5").Value = 11.46645323539191
$ws.Range("M5").Value = 19.42145157915181
$ws.Range("N5").Value = 27.31055495904103
$ws.Range("B6").Value = 21.90323965057105
$ws.Range("C6").Value = 3.48672022416215
$ws.Range("D6").Value = 4.134906936087219
$ws.Range("F6").Value = 56.143063825585
$ws.Range("G6").Value = 3.809741680625869
$ws.Range("J6").Value = 10.7568185773015
$ws.Range("K6").Value = 17.95119977230731
$ws.Range("L6").Value = 11.46737995393564
$ws.Range("M6").Value = 19.42230584121022
$ws.Range("N6").Value = 27.31221883228639
$ws.Range("B7").Value = 21.92415150748381
$ws.Range("C7").Value = 3.554223951726668
$ws.Range("D7").Value = 4.132664356604496
$ws.Range("F7").Value = 56.14364102891841
$ws.Range("G7").Value = 3.808772318838366
$ws.Range("J7").Value = 10.75230329346495
$ws.Range("K7").Value = 17.9645940986396
$ws.Range("L7").Value = 11.46106794108498
$ws.Range("M7").Value = 19.416659161365
$ws.Range("N7").Value = 27.30083432342846
$ws.Range("B8").Value = 22.02755855969775
$ws.Range("C8").Value = 3.836190635682674
$ws.Range("D8").Value = 4.123703963056392
$ws.Range("F8").Value = 56.16573041723102
$ws.Range("G8").Value = 3.804716503503792
$ws.Range("J8").Value = 10.73354634929117
$ws.Range("K8").Value = 18.03263565393399
$ws.Range("L8").Value = 11.43619162099668
$ws.Range("M8").Value = 19.39873530056292
$ws.Range("N8").Value = 27.25466821930199
$ws.Range("B9").Value = 22.26927907917691
$ws.Range("C9").Value = 4.330845605356325
$ws.Range("D9").Value = 4.109476059146802
$ws.Range("F9").Value = 56.27907279887734
$ws.Range("G9").Value = 3.797554234664956
$ws.Range("J9").Value = 10.70093695628312
$ws.Range("K9").Value = 18.19751797773685
$ws.Range("L9").Value = 11.39808295933965
$ws.Range("M9").Value = 19.38866122381246
$ws.Range("N9").Value = 27.17875517643633
$ws.Range("B10").Value = 22.46953693726351
$ws.Range("C10").Value = 4.658119907089354
$ws.Range("D10").Value = 4.10104390811188
$ws.Range("F10").Value = 56.4047800217245
$ws.Range("G10").Value = 3.792768864405269
$ws.Range("J10").Value = 10.67949928332964
$ws.Range("K10").Value = 18.33714805944453
$ws.Range("L10").Value = 11.37656691700152
$ws.Range("M10").Value = 19.39647667911516
$ws.Range("N10").Value = 27.13188146556964
$ws.Range("B11").Value = 22.56530330692735
$ws.Range("C11").Value = 4.799082242313384
$ws.Range("D11").Value = 4.09764315612389
$ws.Range("F11").Value = 56.47110649972041
$ws.Range("G11").Value = 3.790694210176289
$ws.Range("J11").Value = 10.67028905060694
$ws.Range("K11").Value = 18.4045173520458
$ws.Range("L11").Value = 11.36818071659314
$ws.Range("M11").Value = 19.4033167395059
$ws.Range("N11").Value = 27.11248739241948
$ws.Range("B12").Value = 22.60221509328984
$ws.Range("C12").Value = 4.851318959513901
$ws.Range("D12").Value = 4.096417646787603
$ws.Range("F12").Value = 56.49752875901181
$ws.Range("G12").Value = 3.789923203244189
$ws.Range("J12").Value = 10.66687891550706
$ws.Range("K12").Value = 18.43056646535892
$ws.Range("L12").Value = 11.36520611816253
$ws.Range("M12").Value = 19.40637710108
$ws.Range("N12").Value = 27.10542057486501
$ws.Range("B13").Value = 22.59423709393707
$ws.Range("C13").Value = 4.840119718485556
$ws.Range("D13").Value = 4.096678816640289
$ws.Range("F13").Value = 56.49178032401944
$ws.Range("G13").Value = 3.790088604465572
$ws.Range("J13").Value = 10.66760990410779
$ws.Range("K13").Value = 18.42493266578932
$ws.Range("L13").Value = 11.36583781606234
$ws.Range("M13").Value = 19.40569712337438
$ws.Range("N13").Value = 27.10693020904971
$ws.Range("B14").Value = 22.56832723432153
$ws.Range("C14").Value = 4.80340270768501
$ws.Range("D14").Value = 4.097541086127164
$ws.Range("F14").Value = 56.47325415711227
$ws.Range("G14").Value = 3.790630486439484
$ws.Range("J14").Value = 10.67000694385784
$ws.Range("K14").Value = 18.40664974383556
$ws.Range("L14").Value = 11.36793196754111
$ws.Range("M14").Value = 19.40355910753843
$ws.Range("N14").Value = 27.11190044416795
$ws.Range("B15").Value = 22.55254023873944
$ws.Range("C15").Value = 4.780763571883891
$ws.Range("D15").Value = 4.098077353645817
$ws.Range("F15").Value = 56.46207613025911
$ws.Range("G15").Value = 3.790964306284422
$ws.Range("J15").Value = 10.67148529308947
$ws.Range("K15").Value = 18.39552047904853
$ws.Range("L15").Value = 11.36924086721216
$ws.Range("M15").Value = 19.40231067357059
$ws.Range("N15").Value = 27.11498096622363
$ws.Range("B16").Value = 22.46337027749373
$ws.Range("C16").Value = 4.648747982425137
$ws.Range("D16").Value = 4.101274885989739
$ws.Range("F16").Value = 56.40062866483636
$ws.Range("G16").Value = 3.79290649817194
$ws.Range("J16").Value = 10.6801120683857
$ws.Range("K16").Value = 18.33282146613615
$ws.Range("L16").Value = 11.37714313830677
$ws.Range("M16").Value = 19.39609558376099
$ws.Range("N16").Value = 27.13318772256555
$ws.Range("B17").Value = 22.40984620538849
$ws.Range("C17").Value = 4.565729501540309
$ws.Range("D17").Value = 4.103347701147515
$ws.Range("F17").Value = 56.36526818348596
$ws.Range("G17").Value = 3.794124096503371
$ws.Range("J17").Value = 10.68554285893635
$ws.Range("K17").Value = 18.29533290840635
$ws.Range("L17").Value = 11.38234957190259
$ws.Range("M17").Value = 19.39312267671333
$ws.Range("N17").Value = 27.14485097790902
$ws.Range("B18").Value = 22.37950147564801
$ws.Range("C18").Value = 4.517234473525926
$ws.Range("D18").Value = 4.104580897838279
$ws.Range("F18").Value = 56.34579081211035
$ws.Range("G18").Value = 3.794834054516823
$ws.Range("J18").Value = 10.68871752976191
$ws.Range("K18").Value = 18.27413395573264
$ws.Range("L18").Value = 11.38547612529526
$ws.Range("M18").Value = 19.3917220634027
$ws.Range("N18").Value = 27.15174094698009
$ws.Range("B19").Value = 22.36930371273196
$ws.Range("C19").Value = 4.500687096419203
$ws.Range("D19").Value = 4.105005481967098
$ws.Range("F19").Value = 56.33934424102191
$ws.Range("G19").Value = 3.795076090158531
$ws.Range("J19").Value = 10.68980119255437
$ws.Range("K19").Value = 18.26701925216386
$ws.Range("L19").Value = 11.38655739905971
$ws.Range("M19").Value = 19.39130102832155
$ws.Range("N19").Value = 27.15410496038415
$ws.Range("B20").Value = 22.41549847371769
$ws.Range("C20").Value = 4.574644101138114
$ws.Range("D20").Value = 4.103122808432388
$ws.Range("F20").Value = 56.36894331189652
$ws.Range("G20").Value = 3.793993485248094
$ws.Range("J20").Value = 10.68495946371946
$ws.Range("K20").Value = 18.29928612366263
$ws.Range("L20").Value = 11.38178168509346
$ws.Range("M20").Value = 19.39340714870871
$ws.Range("N20").Value = 27.14359061296282
$ws.Range("B21").Value = 22.57592022277272
$ws.Range("C21").Value = 4.814218421111915
$ws.Range("D21").Value = 4.0972861286132
$ws.Range("F21").Value = 56.47866037531979
$ws.Range("G21").Value = 3.790470926444279
$ws.Range("J21").Value = 10.66930077237373
$ws.Range("K21").Value = 18.41200541609706
$ws.Range("L21").Value = 11.36731141125724
$ws.Range("M21").Value = 19.40417435264279
$ws.Range("N21").Value = 27.11043304059817
$ws.Range("B22").Value = 22.68452518033241
$ws.Range("C22").Value = 4.964132929702125
$ws.Range("D22").Value = 4.093834415790802
$ws.Range("F22").Value = 56.55797358968915
$ws.Range("G22").Value = 3.788253907326868
$ws.Range("J22").Value = 10.65951895293819
$ws.Range("K22").Value = 18.48880084643245
$ws.Range("L22").Value = 11.35902602781717
$ws.Range("M22").Value = 19.41395071846771
$ws.Range("N22").Value = 27.09037890793604
$ws.Range("B23").Value = 22.62622487846568
$ws.Range("C23").Value = 4.884731227417599
$ws.Range("D23").Value = 4.095643547847706
$ws.Range("F23").Value = 56.5149497641645
$ws.Range("G23").Value = 3.789429405317659
$ws.Range("J23").Value = 10.66469844404204
$ws.Range("K23").Value = 18.44753305837215
$ws.Range("L23").Value = 11.36334103609303
$ws.Range("M23").Value = 19.40848301789329
$ws.Range("N23").Value = 27.10093431850422
$ws.Range("B24").Value = 22.41294175136615
$ws.Range("C24").Value = 4.570616198184545
$ws.Range("D24").Value = 4.10322435313889
$ws.Range("F24").Value = 56.3672791320874
$ws.Range("G24").Value = 3.794052503612066
$ws.Range("J24").Value = 10.68522305333189
$ws.Range("K24").Value = 18.2974977720101
$ws.Range("L24").Value = 11.38203801144958
$ws.Range("M24").Value = 19.39327757762565
$ws.Range("N24").Value = 27.14415984881433
$ws.Range("B25").Value = 22.19982352323789
$ws.Range("C25").Value = 4.203343371080052
$ws.Range("D25").Value = 4.112968923887691
$ws.Range("F25").Value = 56.24093608031255
$ws.Range("G25").Value = 3.799407695549662
$ws.Range("J25").Value = 10.70931435904219
$ws.Range("K25").Value = 18.14961202049352
$ws.Range("L25").Value = 11.40725212372903
$ws.Range("M25").Value = 19.38870829192001
$ws.Range("N25").Value = 27.19772792669321
